# [ELAB-432] corrections and rewording of some parts^6
#
# Fixes a handful of copy-editing issues in the "Anforderungsliste" sheet:
#  - "Filtierung" -> "Filtrierung" (typo) in three requirement descriptions
#  - "Daten zus Session-Replays" -> "Daten des Session-Replays" (typo)
#  - "3 mal" -> "3-mal" (hyphenation) in the resilience requirement
#
# The shared-string pool is rebuilt/deduplicated automatically by the engine
# on save, so we only need to push the corrected text into the affected
# cells (in the same order the original author touched them) for everything
# else (indices referenced from columns C/D/E/F/G/I) to line up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: "Nutzer ... Logmeldungen ... filtern koennen" description
$ws.Range("G24").Value = "Nutzer des Systems sollen die erfassten Logmeldungen einsehen sowie diese filtern können. Die Filtrierung erfolgt auf Basis der Eigenschaften der Logmeldung (bspw. des Log-Levels)."

# Row 26: "Nutzer ... Fehler ... filtern koennen" description
$ws.Range("G26").Value = "Nutzer des Systems sollen die erfassten Fehler einsehen sowie diese filtern können. Die Filtrierung erfolgt auf Basis der Eigenschaften der Fehler (bspw. der Fehlername)."

# Row 30: Tracingdaten filtering description
$ws.Range("G30").Value = "Die erfassten Tracingdaten sind für die Nutzer des Systems einsehbar und können gefiltert werden. Die Filtrierung erfolgt auf Basis von Eigenschaften der Tracingdaten (wie Name des meldenden Systems)."

# Row 35: Session-Replay partner system description
$ws.Range("G35").Value = "Es existiert ein ""Session-Replay""-Partnersystem, zu dem die Daten des Session-Replays gesendet werden und welches diese analysiert und speichert."

# Row 18: resilience-of-transmission requirement
$ws.Range("G18").Value = "Daten, die der Nachvollziehbarkeit dienen, sollen, wenn möglich, bei einer fehlgeschlagenen Verbindung nicht verworfen werden. Sie sind (mindestens 60s) vorzuhalten und in dieser Zeit sind wiederholt (min. 3-mal) Verbindungsversuche zu unternehmen."

# Restore the selection/viewport the author ended the edit session on.
$ws.Range("G19").Select()
